# Applies the cryptos list update described in the commit message:
# "Updated cryptos list on Mon Jul 17 14:58:02 UTC 2023 with GitHub Actions"
# All Price/Volume cells are plain text in the source data (t="inlineStr"),
# so numeric-looking prices are written with a leading apostrophe (quote-prefix)
# the same way Excel keeps user-typed "numbers" as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.279.39"
$ws.Range("E2").Value = "  -0.26%  "

# Row 3
$ws.Range("D3").Value = "1.913.60"
$ws.Range("E3").Value = "  -1.23%  "

# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").Value = "'0.7406"
$ws.Range("E5").Value = "  -3.10%  "

# Row 6
$ws.Range("D6").Value = "'243.94"
$ws.Range("E6").Value = "  -1.91%  "

# Row 7
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").Value = "'0.3145"
$ws.Range("E8").Value = "  -2.36%  "

# Row 9
$ws.Range("D9").Value = "'27.21"
$ws.Range("E9").Value = "  -3.91%  "

# Row 10
$ws.Range("D10").Value = "'0.07005"
$ws.Range("E10").Value = "  -1.66%  "

# Row 11
$ws.Range("D11").Value = "'0.7836"
$ws.Range("E11").Value = "  -0.89%  "

# Row 12
$ws.Range("D12").Value = "'0.07972"
$ws.Range("E12").Value = "  -0.42%  "

# Row 13
$ws.Range("D13").Value = "1.910.21"
$ws.Range("E13").Value = "  -1.35%  "

# Row 14
$ws.Range("D14").Value = "'5.309"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Value = "'92.07"
$ws.Range("E15").Value = "  -3.02%  "

# Row 16
$ws.Range("D16").Value = "'14.38"
$ws.Range("E16").Value = "  -2.61%  "

# Row 17
$ws.Range("D17").Value = "30.234.28"
$ws.Range("E17").Value = "  -0.38%  "

# Row 18
$ws.Range("D18").Value = "'245.61"
$ws.Range("E18").Value = "  -3.36%  "

# Row 19
$ws.Range("D19").Value = "'5.844"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("D20").Value = "'0.000007866"
$ws.Range("E20").Value = "  -2.21%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  +0.14%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.156.11"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23
$ws.Range("D23").Value = "'0.9995"
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("D24").Value = "'6.658"
$ws.Range("E24").Value = "  -2.59%  "

# Row 25
$ws.Range("D25").Value = "'9.488"
$ws.Range("E25").Value = "  -1.64%  "

# Row 26
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").Value = "'19.01"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("E28").Value = "  -6.00%  "

# Row 29
$ws.Range("D29").Value = "'2.122"
$ws.Range("E29").Value = "  -8.45%  "

# Row 30
$ws.Range("D30").Value = "'1.350"
$ws.Range("E30").Value = "  -0.96%  "

# Row 31
$ws.Range("D31").Value = "'1.546"
$ws.Range("E31").Value = "  +1.10%  "

# Row 32
$ws.Range("D32").Value = "'4.336"
$ws.Range("E32").Value = "  -2.41%  "

# Row 33
$ws.Range("D33").Value = "'4.089"
$ws.Range("E33").Value = "  -1.61%  "

# Row 34
$ws.Range("D34").Value = "'0.05211"
$ws.Range("E34").Value = "  +0.15%  "

# Row 35
$ws.Range("D35").Value = "'1.311"
$ws.Range("E35").Value = "  +1.22%  "

# Row 36
$ws.Range("D36").Value = "'0.7517"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37
$ws.Range("D37").Value = "'2.758"

# Row 38
$ws.Range("D38").Value = "'0.01952"
$ws.Range("E38").Value = "  -1.17%  "

# Row 39
$ws.Range("D39").Value = "'2.798"
$ws.Range("E39").Value = "  -0.32%  "

# Row 40
$ws.Range("E40").Value = "  -1.18%  "

# Row 41
$ws.Range("D41").Value = "'75.99"
$ws.Range("E41").Value = "  -3.26%  "

# Row 42
$ws.Range("D42").Value = "'0.4503"
$ws.Range("E42").Value = "  -0.67%  "

# Row 43
$ws.Range("D43").Value = "'1.946"
$ws.Range("E43").Value = "  -2.75%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("D45").Value = "'7.776"
$ws.Range("E45").Value = "  +2.78%  "

# Row 46
$ws.Range("D46").Value = "'0.8344"
$ws.Range("E46").Value = "  -0.53%  "

# Row 47
$ws.Range("D47").Value = "'9.914"
$ws.Range("E47").Value = "  +0.86%  "

# Row 48
$ws.Range("D48").Value = "'101.20"
$ws.Range("E48").Value = "  -1.17%  "

# Row 49
$ws.Range("D49").Value = "'37.47"

# Row 50
$ws.Range("E50").Value = "  +1.52%  "

# Row 51
$ws.Range("D51").Value = "'939.86"
$ws.Range("E51").Value = "  -5.07%  "

